# adjust the equip addon on levels
# Adds a new "EquipAddon" column (F) to the Exp sheet/table, with header
# rows and 99 data rows of values, and expands the existing table to
# include the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Expand the table ("表1") from A1:E102 to A1:F102 -------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F102"))

# --- Header rows for the new column F ------------------------------------
# Row 1: column header "EquipAddon"
$ws.Range("F1").Value = "EquipAddon"

# Row 2: type row "int" (copy format from E2 which carries style s=1)
$ws.Range("E2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F2").Value = "int"

# Row 3: description row "装备加成" (copy format from E3 which carries style s=2)
$ws.Range("E3").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("F3").Value = "装备加成"

# --- Data rows F4:F102 ----------------------------------------------------
# Values follow the arithmetic progression 3, 8, 13, ... (step 5)
for ($row = 4; $row -le 102; $row++) {
    $ws.Cells.Item($row, 6).Value = 5 * ($row - 4) + 3
}

# --- Restore the selected cell to match the saved view -------------------
$ws.Range("D30").Select()
